$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 30
$ws.Range("G5").Value = 11304.3
$ws.Range("F9").Value = 210
$ws.Range("G9").Value = 6209.7
$ws.Range("B10").Value = 46180.44
$ws.Range("F21").Value = 2
$ws.Range("G21").Value = 90.68000000000001
$ws.Range("B22").Value = 4040.51
$ws.Range("F42").Value = 7
$ws.Range("G42").Value = 156.24
$ws.Range("F43").Value = 1
$ws.Range("G43").Value = 35.32
$ws.Range("B47").Value = 19482.89
$ws.Range("F67").Value = 26
$ws.Range("G67").Value = 1034.8
$ws.Range("B68").Value = 10268.4
$ws.Range("F72").Value = 40
$ws.Range("G72").Value = 8367.200000000001
$ws.Range("F73").Value = 11
$ws.Range("G73").Value = 1045.22
$ws.Range("F78").Value = 43
$ws.Range("G78").Value = 3049.56
$ws.Range("F82").Value = 34
$ws.Range("G82").Value = 1526.6
$ws.Range("F83").Value = 84
$ws.Range("G83").Value = 5610.36
$ws.Range("F84").Value = 35
$ws.Range("G84").Value = 3332.35
$ws.Range("F85").Value = 29
$ws.Range("G85").Value = 659.17
$ws.Range("F86").Value = 52
$ws.Range("G86").Value = 5867.68
$ws.Range("F95").Value = 11
$ws.Range("G95").Value = 4071.98
$ws.Range("F98").Value = 63
$ws.Range("G98").Value = 3584.7
$ws.Range("F102").Value = 34
$ws.Range("G102").Value = 3319.08
$ws.Range("F103").Value = 27
$ws.Range("G103").Value = 1371.6
$ws.Range("F104").Value = 25
$ws.Range("G104").Value = 2561.5
$ws.Range("F109").Value = 166
$ws.Range("G109").Value = 20828.02
$ws.Range("F112").Value = 16
$ws.Range("G112").Value = 756.8
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 0
$ws.Range("B114").Value = 321193.98
$ws.Range("F133").Value = 5
$ws.Range("G133").Value = 211.6
$ws.Range("B138").Value = 1961.1
$ws.Range("B146").Value = 53925
$ws.Range("E146").Value = 79.37
$ws.Range("F146").Value = 1
$ws.Range("G146").Value = 66.44
$ws.Range("B147").Value = 64350
$ws.Range("E147").Value = 70.63
$ws.Range("F147").Value = 2
$ws.Range("G147").Value = 132.88
$ws.Range("B148").Value = 57756
$ws.Range("F148").Value = -100
$ws.Range("G148").Value = -6644
$ws.Range("F174").Value = 31
$ws.Range("G174").Value = 2325.62
$ws.Range("B178").Value = 4958.31
$ws.Range("F186").Value = 2237
$ws.Range("G186").Value = 18902.65
$ws.Range("F187").Value = 1610
$ws.Range("G187").Value = 12863.9
$ws.Range("B189").Value = 40017.17
$ws.Range("F193").Value = 351
$ws.Range("G193").Value = 22744.8
$ws.Range("F194").Value = 89
$ws.Range("G194").Value = 4137.61
$ws.Range("F197").Value = 117
$ws.Range("G197").Value = 5439.33
$ws.Range("F198").Value = 64
$ws.Range("G198").Value = 4808.96
$ws.Range("F199").Value = 3
$ws.Range("G199").Value = 227.91
$ws.Range("B200").Value = 61030.92
$ws.Range("F207").Value = 23
$ws.Range("G207").Value = 6117.08
$ws.Range("F208").Value = 97
$ws.Range("G208").Value = 11119.11
$ws.Range("F213").Value = 20
$ws.Range("G213").Value = 495.2
$ws.Range("B222").Value = 74590.17999999999
$ws.Range("F236").Value = 89
$ws.Range("G236").Value = 3824.33
$ws.Range("F241").Value = 14
$ws.Range("G241").Value = 282.94
$ws.Range("F242").Value = 15
$ws.Range("G242").Value = 481.05
$ws.Range("F245").Value = 5
$ws.Range("G245").Value = 332.95
$ws.Range("B246").Value = 48706
$ws.Range("E246").Value = 39.8
$ws.Range("F246").Value = -144
$ws.Range("G246").Value = -4795.2
$ws.Range("B247").Value = 64973
$ws.Range("E247").Value = 35.4
$ws.Range("F247").Value = 140
$ws.Range("G247").Value = 4662
$ws.Range("F250").Value = 47
$ws.Range("G250").Value = 5384.32
$ws.Range("F256").Value = 13
$ws.Range("G256").Value = 262.73
$ws.Range("F258").Value = 16
$ws.Range("G258").Value = 1040.96
$ws.Range("B274").Value = 122237.14
$ws.Range("F280").Value = 7
$ws.Range("G280").Value = 278.04
$ws.Range("F281").Value = 38
$ws.Range("G281").Value = 6114.2
$ws.Range("F283").Value = 4
$ws.Range("G283").Value = 634.6799999999999
$ws.Range("F290").Value = 90
$ws.Range("G290").Value = 10280.7
$ws.Range("B292").Value = 63520
$ws.Range("E292").Value = 153.4
$ws.Range("F292").Value = 83
$ws.Range("G292").Value = 11975.24
$ws.Range("B293").Value = 55373
$ws.Range("E293").Value = 163.62
$ws.Range("F293").Value = -94
$ws.Range("G293").Value = -13562.32
$ws.Range("B294").Value = 63531
$ws.Range("E294").Value = 152.53
$ws.Range("F294").Value = 80
$ws.Range("G294").Value = 11478.4
$ws.Range("B295").Value = 63571
$ws.Range("F295").Value = 12
$ws.Range("G295").Value = 1721.76
$ws.Range("B296").Value = 57802
$ws.Range("E296").Value = 162.71
$ws.Range("F296").Value = -79
$ws.Range("G296").Value = -11334.92
$ws.Range("F298").Value = 7
$ws.Range("G298").Value = 829.22
$ws.Range("F309").Value = 9
$ws.Range("G309").Value = 1183.41
$ws.Range("B311").Value = 61605
$ws.Range("E311").Value = 133.78
$ws.Range("F311").Value = -13
$ws.Range("G311").Value = -1455.48
$ws.Range("B312").Value = 63563
$ws.Range("E312").Value = 119.04
$ws.Range("F312").Value = 2
$ws.Range("G312").Value = 223.92
$ws.Range("F314").Value = 16
$ws.Range("G314").Value = 2064.16
$ws.Range("F319").Value = 19
$ws.Range("G319").Value = 226.86
$ws.Range("F324").Value = 80
$ws.Range("G324").Value = 4730.4
$ws.Range("F325").Value = 197
$ws.Range("G325").Value = 27150.54
$ws.Range("F328").Value = 1680
$ws.Range("G328").Value = 35330.4
$ws.Range("F329").Value = 66
$ws.Range("G329").Value = 10626
$ws.Range("F330").Value = 13
$ws.Range("G330").Value = 6837.35
$ws.Range("F333").Value = 803
$ws.Range("G333").Value = 137577.99
$ws.Range("F335").Value = 14
$ws.Range("G335").Value = 2603.86
$ws.Range("B339").Value = 418027.06
$ws.Range("F343").Value = 12
$ws.Range("G343").Value = 1940.76
$ws.Range("B346").Value = 15968.62
$ws.Range("F363").Value = 60
$ws.Range("G363").Value = 1275
$ws.Range("F364").Value = 29
$ws.Range("G364").Value = 3475.07
$ws.Range("F366").Value = 43
$ws.Range("G366").Value = 5904.76
$ws.Range("F370").Value = 63
$ws.Range("G370").Value = 3383.1
$ws.Range("F372").Value = 29
$ws.Range("G372").Value = 1564.55
$ws.Range("F376").Value = 96
$ws.Range("G376").Value = 672
$ws.Range("F377").Value = 14
$ws.Range("G377").Value = 3661.14
$ws.Range("F378").Value = 17
$ws.Range("G378").Value = 1260.72
$ws.Range("F382").Value = 183
$ws.Range("G382").Value = 7870.83
$ws.Range("F383").Value = 99
$ws.Range("G383").Value = 8243.73
$ws.Range("F385").Value = 121
$ws.Range("G385").Value = 8634.559999999999
$ws.Range("F389").Value = 7
$ws.Range("G389").Value = 391.02
$ws.Range("F394").Value = 103
$ws.Range("G394").Value = 21721.67
$ws.Range("B395").Value = 282523.17
$ws.Range("F398").Value = 101
$ws.Range("G398").Value = 2132.11
$ws.Range("F399").Value = 33
$ws.Range("G399").Value = 5086.95
$ws.Range("B402").Value = 10046.49
$ws.Range("F405").Value = 64
$ws.Range("G405").Value = 517.76
$ws.Range("F408").Value = 13
$ws.Range("G408").Value = 416.26
$ws.Range("B411").Value = 4249.86
$ws.Range("F422").Value = 14
$ws.Range("G422").Value = 3707.62
$ws.Range("F426").Value = 246
$ws.Range("G426").Value = 23763.6
$ws.Range("B430").Value = 64771.31
$ws.Range("F436").Value = 10
$ws.Range("G436").Value = 1811.6
$ws.Range("F439").Value = 113
$ws.Range("G439").Value = 2678.1
$ws.Range("F442").Value = 71
$ws.Range("G442").Value = 3617.45
$ws.Range("F443").Value = 161
$ws.Range("G443").Value = 5098.87
$ws.Range("B448").Value = 49249.42
$ws.Range("F465").Value = 23
$ws.Range("G465").Value = 309.35
$ws.Range("F470").Value = 208
$ws.Range("G470").Value = 2664.48
$ws.Range("B472").Value = 45695
$ws.Range("E472").Value = 23.58
$ws.Range("F472").Value = -36
$ws.Range("G472").Value = -710.28
$ws.Range("B473").Value = 64915
$ws.Range("E473").Value = 20.98
$ws.Range("F473").Value = 0
$ws.Range("G473").Value = 0
$ws.Range("F475").Value = 22
$ws.Range("G475").Value = 281.82
$ws.Range("F477").Value = 157
$ws.Range("G477").Value = 3097.61
$ws.Range("B479").Value = 64927
$ws.Range("E479").Value = 17.26
$ws.Range("F479").Value = 256
$ws.Range("G479").Value = 4152.32
$ws.Range("B480").Value = 45718
$ws.Range("E480").Value = 19.38
$ws.Range("F480").Value = -294
$ws.Range("G480").Value = -4768.68
$ws.Range("B485").Value = 45709
$ws.Range("E485").Value = 15.69
$ws.Range("F485").Value = -300
$ws.Range("G485").Value = -3945
$ws.Range("B486").Value = 64925
$ws.Range("E486").Value = 13.97
$ws.Range("F486").Value = 253
$ws.Range("G486").Value = 3326.95
$ws.Range("B487").Value = 45702
$ws.Range("E487").Value = 31.43
$ws.Range("F487").Value = -215
$ws.Range("G487").Value = -5654.5
$ws.Range("B488").Value = 64919
$ws.Range("E488").Value = 27.97
$ws.Range("F488").Value = 157
$ws.Range("G488").Value = 4129.1
$ws.Range("F489").Value = 96
$ws.Range("G489").Value = 1577.28
$ws.Range("F490").Value = 271
$ws.Range("G490").Value = 3991.83
$ws.Range("B492").Value = 7047.29
$ws.Range("F494").Value = 17
$ws.Range("G494").Value = 861.39
$ws.Range("F495").Value = 53
$ws.Range("G495").Value = 1658.9
$ws.Range("F501").Value = 16
$ws.Range("G501").Value = 2549.12
$ws.Range("F502").Value = 13
$ws.Range("G502").Value = 2071.16
$ws.Range("B508").Value = 25910.94
$ws.Range("F510").Value = 15
$ws.Range("G510").Value = 209.7
$ws.Range("F511").Value = 119
$ws.Range("G511").Value = 2885.75
$ws.Range("F516").Value = 33
$ws.Range("G516").Value = 533.9400000000001
$ws.Range("F526").Value = 21
$ws.Range("G526").Value = 2210.67
$ws.Range("B528").Value = 25213.31
$ws.Range("F547").Value = 175
$ws.Range("G547").Value = 3473.75
$ws.Range("F549").Value = 7
$ws.Range("G549").Value = 115.71
$ws.Range("F551").Value = 72
$ws.Range("G551").Value = 2350.08
$ws.Range("B552").Value = 10486.77
$ws.Range("F569").Value = 4
$ws.Range("G569").Value = 83.68000000000001
$ws.Range("B582").Value = 24498.69
$ws.Range("B585").Value = 60025
$ws.Range("E585").Value = 37.22
$ws.Range("F585").Value = -98
$ws.Range("G585").Value = -3217.34
$ws.Range("B586").Value = 64833
$ws.Range("E586").Value = 34.9
$ws.Range("F586").Value = 96
$ws.Range("G586").Value = 3151.68
$ws.Range("F607").Value = 11
$ws.Range("G607").Value = 4690.18
$ws.Range("F612").Value = 94
$ws.Range("G612").Value = 8303.959999999999
$ws.Range("B614").Value = 48252.31
$ws.Range("F627").Value = 24
$ws.Range("G627").Value = 1036.32
$ws.Range("B629").Value = 2693.02
$ws.Range("F636").Value = 0
$ws.Range("G636").Value = 0
$ws.Range("B638").Value = 75566.64
$ws.Range("F650").Value = 403
$ws.Range("G650").Value = 32393.14
$ws.Range("B651").Value = 42539.83
$ws.Range("F677").Value = 0
$ws.Range("G677").Value = 0
$ws.Range("B679").Value = 41423.91
$ws.Range("F682").Value = 7
$ws.Range("G682").Value = 2236.64
$ws.Range("F686").Value = 20
$ws.Range("G686").Value = 6908.6
$ws.Range("B696").Value = 43197.68
$ws.Range("F698").Value = 108
$ws.Range("G698").Value = 8808.48
$ws.Range("F701").Value = 252
$ws.Range("G701").Value = 36068.76
$ws.Range("F702").Value = 100
$ws.Range("G702").Value = 8156
$ws.Range("F703").Value = 104
$ws.Range("G703").Value = 10586.16
$ws.Range("F704").Value = 40
$ws.Range("G704").Value = 5324
$ws.Range("F705").Value = 130
$ws.Range("G705").Value = 9838.4
$ws.Range("F707").Value = 198
$ws.Range("G707").Value = 4300.56
$ws.Range("F708").Value = 77
$ws.Range("G708").Value = 2870.56
$ws.Range("F710").Value = 107
$ws.Range("G710").Value = 7442.92
$ws.Range("F713").Value = 581
$ws.Range("G713").Value = 78440.81
$ws.Range("F715").Value = 509
$ws.Range("G715").Value = 61441.39
$ws.Range("B716").Value = 248061.38
$ws.Range("F720").Value = 88
$ws.Range("G720").Value = 14422.32
$ws.Range("F722").Value = 15
$ws.Range("G722").Value = 1632.15
$ws.Range("F724").Value = 73
$ws.Range("G724").Value = 10979.93
$ws.Range("F725").Value = 33
$ws.Range("G725").Value = 4828.89
$ws.Range("F731").Value = 0
$ws.Range("G731").Value = 0
$ws.Range("F734").Value = 10
$ws.Range("G734").Value = 472.1
$ws.Range("F735").Value = 170
$ws.Range("G735").Value = 5620.2
$ws.Range("F738").Value = 186
$ws.Range("G738").Value = 9279.540000000001
$ws.Range("F740").Value = 122
$ws.Range("G740").Value = 6952.78
$ws.Range("F741").Value = 161
$ws.Range("G741").Value = 38929.8
$ws.Range("F742").Value = 143
$ws.Range("G742").Value = 8149.57
$ws.Range("B743").Value = 131715.34
$ws.Range("F768").Value = 3771
$ws.Range("G768").Value = 615087.8100000001
$ws.Range("F769").Value = 98
$ws.Range("G769").Value = 17264.66
$ws.Range("F770").Value = 646
$ws.Range("G770").Value = 182734.02
$ws.Range("F771").Value = 546
$ws.Range("G771").Value = 78978.89999999999
$ws.Range("F772").Value = 10
$ws.Range("G772").Value = 771.4
$ws.Range("F773").Value = 3
$ws.Range("G773").Value = 443.58
$ws.Range("B775").Value = 927757.9
$ws.Range("F778").Value = 161
$ws.Range("G778").Value = 23507.61
$ws.Range("F784").Value = 11
$ws.Range("G784").Value = 1524.05
$ws.Range("F787").Value = 102
$ws.Range("G787").Value = 12785.7
$ws.Range("B792").Value = 99715.37
$ws.Range("B793").Value = 3897375.51
$ws.Range("B794").Value = 3897375.51
